# Aggiunta slide valutazione risultati
#
# The new slide ("Valutazione dei risultati ottenuti") reuses the same
# title-slide layout / formatting as the preceding "Esempi applicativi"
# slide (slide 8), so we duplicate that slide (placing the copy right
# after it, i.e. at the end of the deck) and just swap the caption text.

$p = $ppt.ActivePresentation

# Slide 8 = "Esempi applicativi" (Sottotitolo 2 / Titolo 1 / CasellaDiTesto 1)
$srcSlide = $p.Slides.Item(8)

$dup = $srcSlide.Duplicate()
$newSlide = $dup.Item(1)

# Locate the free-floating caption textbox (not a placeholder) and update
# its text + name; placeholders ("Sottotitolo 2" / "Titolo 1") are left as
# duplicated, since their content/formatting is unchanged on the new slide.
for ($i = 1; $i -le $newSlide.Shapes.Count; $i++) {
    $shp = $newSlide.Shapes.Item($i)
    if ($shp.Name -like "CasellaDiTesto*") {
        $shp.Name = "CasellaDiTesto 4"
        $shp.TextFrame.TextRange.Text = "Valutazione dei risultati ottenuti"
    }
}
